$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.986.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.743.35'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.66%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4990'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3576'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.47'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07265'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.063'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9996'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.29'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("E14").Value = '  -0.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.743.37'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.871'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.51%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.99%  '
$ws.Range("E18").Value = '  -2.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06404'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9995'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.55%  '
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.733'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.060.97'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.048'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.70%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.943.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.153'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.29'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.060'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09500'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.38%  '
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.401'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02205'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05900'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.14'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2000'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.421'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.765'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6036'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9992'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.110'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.501'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.91'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.604'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5657'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.81'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.862'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.105'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06680'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.07%  '
